$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Slit2"
$ws.Cells.Item(2, 3).Value2 = "Gpc1"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 2
$ws.Cells.Item(2, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(2, 7).Value2 = 0.01246433333333333
$ws.Cells.Item(2, 8).Value2 = 0.037393
$ws.Cells.Item(2, 9).Value2 = 0.0065371131913745
$ws.Cells.Item(2, 10).Value2 = 0.006537113191374499
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 2.127396333333333
$ws.Cells.Item(2, 14).Value2 = 6.382189
$ws.Cells.Item(2, 15).Value2 = 0.06137654768277986
$ws.Cells.Item(2, 16).Value2 = 0.06137654768277986
$ws.Cells.Item(2, 17).Value2 = 0.02651657703077778
$ws.Cells.Item(2, 18).Value2 = 0.238649193277
$ws.Cells.Item(2, 19).Value2 = 0.0004012254394981263
$ws.Cells.Item(2, 20).Value2 = 0.0004012254394981262

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Slit2"
$ws.Cells.Item(3, 3).Value2 = "Gpc1"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 7).Value2 = 0.01246433333333333
$ws.Cells.Item(3, 8).Value2 = 0.037393
$ws.Cells.Item(3, 9).Value2 = 0.0065371131913745
$ws.Cells.Item(3, 10).Value2 = 0.006537113191374499
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 3.721182333333333
$ws.Cells.Item(3, 14).Value2 = 11.163547
$ws.Cells.Item(3, 15).Value2 = 0.1073581454191429
$ws.Cells.Item(3, 16).Value2 = 0.1073581454191429
$ws.Cells.Item(3, 17).Value2 = 0.04638205699677777
$ws.Cells.Item(3, 18).Value2 = 0.4174385129709999
$ws.Cells.Item(3, 19).Value2 = 0.0007018123486209807
$ws.Cells.Item(3, 20).Value2 = 0.0007018123486209806

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Slit2"
$ws.Cells.Item(4, 3).Value2 = "Gpc1"
$ws.Cells.Item(4, 4).Value2 = "sCs"
$ws.Cells.Item(4, 5).Value2 = 2
$ws.Cells.Item(4, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 7).Value2 = 0.01246433333333333
$ws.Cells.Item(4, 8).Value2 = 0.037393
$ws.Cells.Item(4, 9).Value2 = 0.0065371131913745
$ws.Cells.Item(4, 10).Value2 = 0.006537113191374499
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 28.81280933333333
$ws.Cells.Item(4, 14).Value2 = 86.43842799999999
$ws.Cells.Item(4, 15).Value2 = 0.8312653068980773
$ws.Cells.Item(4, 16).Value2 = 0.8312653068980772
$ws.Cells.Item(4, 17).Value2 = 0.3591324598004444
$ws.Cells.Item(4, 18).Value2 = 3.232192138203999
$ws.Cells.Item(4, 19).Value2 = 0.005434075403255394
$ws.Cells.Item(4, 20).Value2 = 0.005434075403255392

# Row 5
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Slit2"
$ws.Cells.Item(5, 3).Value2 = "Gpc1"
$ws.Cells.Item(5, 4).Value2 = "ECs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 0.9943730000000001
$ws.Cells.Item(5, 8).Value2 = 2.983119
$ws.Cells.Item(5, 9).Value2 = 0.5215143627507798
$ws.Cells.Item(5, 10).Value2 = 0.5215143627507798
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 2.127396333333333
$ws.Cells.Item(5, 14).Value2 = 6.382189
$ws.Cells.Item(5, 15).Value2 = 0.06137654768277986
$ws.Cells.Item(5, 16).Value2 = 0.06137654768277986
$ws.Cells.Item(5, 17).Value2 = 2.115425474165667
$ws.Cells.Item(5, 18).Value2 = 19.038829267491
$ws.Cells.Item(5, 19).Value2 = 0.03200875115262779
$ws.Cells.Item(5, 20).Value2 = 0.03200875115262779

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Slit2"
$ws.Cells.Item(6, 3).Value2 = "Gpc1"
$ws.Cells.Item(6, 4).Value2 = "FAPs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 0.9943730000000001
$ws.Cells.Item(6, 8).Value2 = 2.983119
$ws.Cells.Item(6, 9).Value2 = 0.5215143627507798
$ws.Cells.Item(6, 10).Value2 = 0.5215143627507798
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 3.721182333333333
$ws.Cells.Item(6, 14).Value2 = 11.163547
$ws.Cells.Item(6, 15).Value2 = 0.1073581454191429
$ws.Cells.Item(6, 16).Value2 = 0.1073581454191429
$ws.Cells.Item(6, 17).Value2 = 3.700243240343667
$ws.Cells.Item(6, 18).Value2 = 33.302189163093
$ws.Cells.Item(6, 19).Value2 = 0.05598881479436985
$ws.Cells.Item(6, 20).Value2 = 0.05598881479436984

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Slit2"
$ws.Cells.Item(7, 3).Value2 = "Gpc1"
$ws.Cells.Item(7, 4).Value2 = "sCs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 0.9943730000000001
$ws.Cells.Item(7, 8).Value2 = 2.983119
$ws.Cells.Item(7, 9).Value2 = 0.5215143627507798
$ws.Cells.Item(7, 10).Value2 = 0.5215143627507798
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 28.81280933333333
$ws.Cells.Item(7, 14).Value2 = 86.43842799999999
$ws.Cells.Item(7, 15).Value2 = 0.8312653068980773
$ws.Cells.Item(7, 16).Value2 = 0.8312653068980772
$ws.Cells.Item(7, 17).Value2 = 28.65067965521467
$ws.Cells.Item(7, 18).Value2 = 257.856116896932
$ws.Cells.Item(7, 19).Value2 = 0.4335167968037822
$ws.Cells.Item(7, 20).Value2 = 0.4335167968037821

# Row 8
$ws.Cells.Item(8, 1).Value2 = "sCs"
$ws.Cells.Item(8, 2).Value2 = "Slit2"
$ws.Cells.Item(8, 3).Value2 = "Gpc1"
$ws.Cells.Item(8, 4).Value2 = "ECs"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 0.8998656666666666
$ws.Cells.Item(8, 8).Value2 = 2.699597
$ws.Cells.Item(8, 9).Value2 = 0.4719485240578458
$ws.Cells.Item(8, 10).Value2 = 0.4719485240578457
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 2.127396333333333
$ws.Cells.Item(8, 14).Value2 = 6.382189
$ws.Cells.Item(8, 15).Value2 = 0.06137654768277986
$ws.Cells.Item(8, 16).Value2 = 0.06137654768277986
$ws.Cells.Item(8, 17).Value2 = 1.914370919759222
$ws.Cells.Item(8, 18).Value2 = 17.229338277833
$ws.Cells.Item(8, 19).Value2 = 0.02896657109065395
$ws.Cells.Item(8, 20).Value2 = 0.02896657109065394

# Row 9
$ws.Cells.Item(9, 1).Value2 = "sCs"
$ws.Cells.Item(9, 2).Value2 = "Slit2"
$ws.Cells.Item(9, 3).Value2 = "Gpc1"
$ws.Cells.Item(9, 4).Value2 = "FAPs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 0.8998656666666666
$ws.Cells.Item(9, 8).Value2 = 2.699597
$ws.Cells.Item(9, 9).Value2 = 0.4719485240578458
$ws.Cells.Item(9, 10).Value2 = 0.4719485240578457
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 3.721182333333333
$ws.Cells.Item(9, 14).Value2 = 11.163547
$ws.Cells.Item(9, 15).Value2 = 0.1073581454191429
$ws.Cells.Item(9, 16).Value2 = 0.1073581454191429
$ws.Cells.Item(9, 17).Value2 = 3.348564221173222
$ws.Cells.Item(9, 18).Value2 = 30.137077990559
$ws.Cells.Item(9, 19).Value2 = 0.05066751827615205
$ws.Cells.Item(9, 20).Value2 = 0.05066751827615204

# Row 10
$ws.Cells.Item(10, 1).Value2 = "sCs"
$ws.Cells.Item(10, 2).Value2 = "Slit2"
$ws.Cells.Item(10, 3).Value2 = "Gpc1"
$ws.Cells.Item(10, 4).Value2 = "sCs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 0.8998656666666666
$ws.Cells.Item(10, 8).Value2 = 2.699597
$ws.Cells.Item(10, 9).Value2 = 0.4719485240578458
$ws.Cells.Item(10, 10).Value2 = 0.4719485240578457
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 28.81280933333333
$ws.Cells.Item(10, 14).Value2 = 86.43842799999999
$ws.Cells.Item(10, 15).Value2 = 0.8312653068980773
$ws.Cells.Item(10, 16).Value2 = 0.8312653068980772
$ws.Cells.Item(10, 17).Value2 = 25.92765787927955
$ws.Cells.Item(10, 18).Value2 = 233.3489209135159
$ws.Cells.Item(10, 19).Value2 = 0.3923144346910398
$ws.Cells.Item(10, 20).Value2 = 0.3923144346910397

